$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the full data range so numeric-looking strings
# (e.g. "252.90", "0.06840") are preserved exactly, matching the source feed formatting.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.109.26'
$ws.Range('E2').Value = '  +4.18%  '
$ws.Range('D3').Value = '1.906.77'
$ws.Range('E3').Value = '  +5.20%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '252.90'
$ws.Range('E5').Value = '  +1.95%  '
$ws.Range('D6').Value = '0.9994'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.5098'
$ws.Range('E7').Value = '  +3.09%  '
$ws.Range('D8').Value = '45.18'
$ws.Range('E8').Value = '  +4.37%  '
$ws.Range('D9').Value = '0.3022'
$ws.Range('E9').Value = '  +8.47%  '
$ws.Range('D10').Value = '0.06840'
$ws.Range('E10').Value = '  +6.39%  '
$ws.Range('D11').Value = '1.906.60'
$ws.Range('E11').Value = '  +5.50%  '
$ws.Range('D12').Value = '17.29'
$ws.Range('E12').Value = '  +2.77%  '
$ws.Range('E13').Value = '  +3.31%  '
$ws.Range('D14').Value = '0.6915'
$ws.Range('E14').Value = '  +6.87%  '
$ws.Range('D15').Value = '86.99'
$ws.Range('E15').Value = '  +3.49%  '
$ws.Range('D16').Value = '4.928'
$ws.Range('E16').Value = '  +5.06%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.000008372'
$ws.Range('E17').Value = '  +13.89%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '30.112.78'
$ws.Range('E18').Value = '  +4.15%  '
$ws.Range('D19').Value = '0.9998'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('E20').Value = '  +6.33%  '
$ws.Range('D21').Value = '2.153.75'
$ws.Range('E21').Value = '  +5.32%  '
$ws.Range('D22').Value = '0.9986'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '4.822'
$ws.Range('E23').Value = '  +5.12%  '
$ws.Range('D24').Value = '5.745'
$ws.Range('D25').Value = '9.294'
$ws.Range('E25').Value = '  +5.03%  '
$ws.Range('D26').Value = '147.51'
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').Value = '134.61'
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('D28').Value = '17.14'
$ws.Range('E28').Value = '  +4.33%  '
$ws.Range('D29').Value = '2.005'
$ws.Range('E29').Value = '  +5.58%  '
$ws.Range('D30').Value = '1.403'
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('D31').Value = '4.287'
$ws.Range('E31').Value = '  +3.32%  '
$ws.Range('D32').Value = '0.08861'
$ws.Range('E32').Value = '  +6.05%  '
$ws.Range('D33').Value = '4.006'
$ws.Range('E33').Value = '  +4.77%  '
$ws.Range('D34').Value = '0.05060'
$ws.Range('E34').Value = '  +1.88%  '
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('D36').Value = '0.7238'
$ws.Range('E36').Value = '  +7.17%  '
$ws.Range('D37').Value = '2.691'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('D38').Value = '2.823'
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('D39').Value = '2.276'
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('D40').Value = '0.9617'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D41').Value = '0.01694'
$ws.Range('E41').Value = '  +6.20%  '
$ws.Range('D42').Value = '6.103'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').Value = '0.4314'
$ws.Range('E43').Value = '  +5.09%  '
$ws.Range('D44').Value = '104.73'
$ws.Range('E44').Value = '  +4.79%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '7.684'
$ws.Range('E46').Value = '  +7.12%  '
$ws.Range('D48').Value = '0.05751'
$ws.Range('E48').Value = '  +4.16%  '
$ws.Range('E49').Value = '  +4.57%  '
$ws.Range('D50').Value = '8.451'
$ws.Range('E50').Value = '  +3.78%  '
$ws.Range('D51').Value = '0.3824'
$ws.Range('E51').Value = '  +4.98%  '

# Restore default (General) styling so cells match the original workbook look.
$ws.Range("B2:E51").Style = "Normal"
